$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free direct cell updates matching the diff.
# For numeric-looking text values (Price column), force text storage
# so Excel does not silently coerce them into floating point numbers
# (which would corrupt the exact decimal text, e.g. "589.61" -> 589.6100000000001).

$ws.Range("D2").Value = "64.340.57"
$ws.Range("D3").Value = "3.501.47"
$ws.Range("E3").Value = "  +0.34%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "589.61"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "134.09"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.10%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "7.72"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +7.06%  "
$ws.Range("E10").Value = "  +0.33%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.388"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +2.75%  "
$ws.Range("D12").Value = "4.095.91"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "3.501.07"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").Value = "64.286.76"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "25.31"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.22%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "10.04"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("E19").Value = "  +0.32%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "13.55"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.60%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "386.65"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.13%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.579"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +2.51%  "
$ws.Range("D23").Value = "3.639.33"
$ws.Range("E23").Value = "  +0.18%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "74.29"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  +0.08%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "5.73"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("E27").Value = "  +1.92%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.03%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "7.29"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("E30").Value = "  +1.06%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.49"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.07%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "8.13"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("E33").Value = "  +4.16%  "
$ws.Range("D34").Value = "3.526.98"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("E36").Value = "  -0.56%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "5.37"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +2.49%  "
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("E39").Value = "  +0.41%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "164.65"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("E41").Value = "  +0.57%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.807"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("E43").Value = "  +0.01%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "4.41"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.18"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "24.36"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -4.76%  "
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").Value = "2.425.44"
$ws.Range("E48").Value = "  -2.11%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "6.81"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.99%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.919"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +1.62%  "
$ws.Range("E51").Value = "  -0.21%  "
